# Generate Report for Handback
# Updates the timestamps / status recorded in the handback-status report.
# Each "duplicate" data row (3 and 5) on a sheet shares the same
# shared-string text, so both rows must be updated together to keep the
# in-memory shared-string table consistent with the sheet data.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" column (G) for the
# 7782a12c-e33d-42d3-ab9d-b4bafe1ecf16.md row and its duplicate row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 10:18:36"
$wsOverview.Range("G5").Value = "2016-09-01 10:18:36"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): "ht" -> "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H3").Value = "2016-09-01 10:18:31"
$wsZhCn.Range("H5").Value = "2016-09-01 10:18:31"
# Correspond Handback DateTime column (K)
$wsZhCn.Range("K3").Value = "2016-09-01 10:18:48"
$wsZhCn.Range("K5").Value = "2016-09-01 10:18:48"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column (E): "ht" -> "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsDeDe.Range("H3").Value = "2016-09-01 10:18:36"
$wsDeDe.Range("H5").Value = "2016-09-01 10:18:36"
# Correspond Handback DateTime column (K)
$wsDeDe.Range("K3").Value = "2016-09-01 10:18:55"
$wsDeDe.Range("K5").Value = "2016-09-01 10:18:55"
